$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update cell values (order matters for shared-string table layout):
# B8 email changes from "sals1234@gmail.com" to "emailsalsaaaa@gmail.com"
# B7 password changes from "1234sals!s" to "emailsalsa2"
$ws.Range("B8").Value = "emailsalsaaaa@gmail.com"
$ws.Range("B7").Value = "emailsalsa2"

# Update the selection on the sheet to B7
$ws.Range("B7").Select()

# Update the saved window position
$excel.ActiveWindow.Left = 4950
$excel.ActiveWindow.Top = 3315
